$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Content.Find.Execute("2025-07-11 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-07-12 Saturday", 2)

# Update the table of division problems.
# The table has 20 rows x 5 columns; only rows 1, 5, 9, 13, 17 contain data.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "159÷6=26, 3"
$t.Cell(1, 2).Range.Text = "498÷7=71, 1"
$t.Cell(1, 3).Range.Text = "150÷9=16, 6"
$t.Cell(1, 4).Range.Text = "823÷7=117, 4"
$t.Cell(1, 5).Range.Text = "923÷5=184, 3"

$t.Cell(5, 1).Range.Text = "260÷2=130, 0"
$t.Cell(5, 2).Range.Text = "246÷7=35, 1"
$t.Cell(5, 3).Range.Text = "651÷7=93, 0"
$t.Cell(5, 4).Range.Text = "724÷4=181, 0"
$t.Cell(5, 5).Range.Text = "256÷4=64, 0"

$t.Cell(9, 1).Range.Text = "116÷6=19, 2"
$t.Cell(9, 2).Range.Text = "980÷3=326, 2"
$t.Cell(9, 3).Range.Text = "892÷5=178, 2"
$t.Cell(9, 4).Range.Text = "668÷4=167, 0"
$t.Cell(9, 5).Range.Text = "746÷7=106, 4"

$t.Cell(13, 1).Range.Text = "843÷4=210, 3"
$t.Cell(13, 2).Range.Text = "825÷7=117, 6"
$t.Cell(13, 3).Range.Text = "360÷4=90, 0"
$t.Cell(13, 4).Range.Text = "661÷4=165, 1"
$t.Cell(13, 5).Range.Text = "526÷5=105, 1"

$t.Cell(17, 1).Range.Text = "791÷8=98, 7"
$t.Cell(17, 2).Range.Text = "278÷2=139, 0"
$t.Cell(17, 3).Range.Text = "262÷6=43, 4"
$t.Cell(17, 4).Range.Text = "439÷7=62, 5"
$t.Cell(17, 5).Range.Text = "107÷6=17, 5"
